$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 43 (A:F) into a brand-new row 44 by copy + insert, which
# clones cell formatting (style indices) and adjusts the relative formula
# in column D, instead of recomputing styles from scratch (which would
# otherwise mint new style entries in xl/styles.xml).
$ws.Range("A43:F43").Copy()
$ws.Range("A44:F44").Insert(-4121)

# Overwrite the copied values with the new entry's actual data:
# 4.4.2020, 11:45 - 12:00, Task "Mem Ctrl 2" (new), Unit "Arch and TB".
$ws.Range("B44").Value2 = 0.48958333333333331
$ws.Range("C44").Value2 = 0.5
$ws.Range("D44").Formula = "=C44-B44"
$ws.Range("E44").Value2 = "Mem Ctrl 2"
$ws.Range("F44").Value2 = "Arch and TB"

$ws.Range("F44").Select()
